# Security Program Tracking Template.xlsx - apply commit changes
# - Rename "ERP System Logs (SAP)" entry on the "Log Collection" sheet to the
#   fuller "Enterprise Resource Planning (ERP) System Logs (SAP)" name and
#   widen column A to fit it.
# - Populate the previously-empty "SIEM" sheet with the SIEM capability
#   tracking table (header + 7 rows) and appropriate column widths/styles.
# - Update sheet selections and make "SIEM" the active tab, matching the
#   saved workbook view state.

$wb = $excel.ActiveWorkbook

$wsLog  = $wb.Worksheets("Log Collection")
$wsSiem = $wb.Worksheets("SIEM")

# ---------------------------------------------------------------------
# 1. "Log Collection" sheet: rename the ERP row, widen column A
# ---------------------------------------------------------------------
$wsLog.Range("A25").Value = "Enterprise Resource Planning (ERP) System Logs (SAP)"

# stored OOXML column width = ColumnWidth + 5/6, so back-solve for the
# target stored width of 45.1640625
$wsLog.Columns.Item(1).ColumnWidth = 45.1640625 - 0.8333333333333334

# ---------------------------------------------------------------------
# 2. "SIEM" sheet: build out the capability tracking table
# ---------------------------------------------------------------------
$siemHeader = @('Capability', 'Description', 'Availability', 'Maturity Level', 'Owner', 'Improvement Plan', 'Use Cases')
$siemRows = @(
  @('Detection Capabilities', 'Ability to detect various threats, such as malware, anomalies, and insider threats using the SIEM.', 'Partially Available', 'Stage 2 - Procedural', 'SIEM Team', 'Improve detection rules for insider threats', 'Threat Detection, Compliance'),
  @('Threat Hunting', 'Ability to proactively search for indicators of compromise (IOCs) and unknown threats.', 'Available', 'Stage 3 - Innovative', 'Threat Intel Team', 'Increase frequency of threat hunts', 'Threat Detection, Incident Response'),
  @('Custom Dashboards', 'Customizable dashboards to monitor specific activities or use cases.', 'Available', 'Stage 3 - Innovative', 'SIEM Admin Team', 'Develop more user-specific dashboards', 'Visibility, Reporting, Compliance'),
  @('Alert Creation', 'Ability to create and manage alerts for suspicious activities detected in logs.', 'Available', 'Stage 3 - Innovative', 'SOC Team', 'Refine alert thresholds to reduce false positives', 'Threat Detection, Incident Response'),
  @('Admin Functions', 'Administrative tasks such as user management, log retention, and configuration changes.', 'Available', 'Stage 3 - Innovative', 'SIEM Admin Team', 'Streamline user role assignment', 'System Management, Compliance'),
  @('Upkeep', 'Regular maintenance, software updates, and system health checks for the SIEM.', 'Available', 'Stage 2 - Procedural', 'SIEM Maintenance Team', 'Automate routine maintenance tasks', 'System Stability, Compliance'),
  @('Threat Intel Feed Search', 'Capability to search for and correlate threat intelligence feeds with internal activity to identify threats.', 'Available', 'Stage 3 - Innovative', 'Threat Intel Team', 'Integrate additional threat feeds', 'Threat Detection, Threat Intelligence')
)

# Pull matching header / body styles (bold+bordered vs plain+bordered) from
# the already-formatted "Log Collection" sheet so the new table looks the
# same as the rest of the workbook.
$wsLog.Range("A1").Copy()
$wsSiem.Range("A1:G1").PasteSpecial(-4122)

$wsLog.Range("A2").Copy()
$wsSiem.Range("A2:G8").PasteSpecial(-4122)

for ($c = 1; $c -le $siemHeader.Length; $c++) {
    $wsSiem.Cells.Item(1, $c).Value = $siemHeader[$c - 1]
}

for ($r = 0; $r -lt $siemRows.Length; $r++) {
    $row = $siemRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsSiem.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Column widths (stored width = ColumnWidth + 5/6)
$pad = 0.8333333333333334
$wsSiem.Columns.Item(1).ColumnWidth = 20.83203125 - $pad
$wsSiem.Columns.Item(2).ColumnWidth = 85 - $pad
$wsSiem.Columns.Item(3).ColumnWidth = 15.5 - $pad
$wsSiem.Columns.Item(4).ColumnWidth = 17.33203125 - $pad
$wsSiem.Columns.Item(5).ColumnWidth = 21 - $pad
$wsSiem.Columns.Item(6).ColumnWidth = 39.83203125 - $pad
$wsSiem.Columns.Item(7).ColumnWidth = 31.1640625 - $pad

# ---------------------------------------------------------------------
# 3. Selections / active tab: "Log Collection" selects D31, "SIEM" selects
#    B16 and becomes the active sheet/tab (matches saved workbookView).
# ---------------------------------------------------------------------
$wsLog.Range("D31").Select()
$wsSiem.Range("B16").Select()
$wsSiem.Activate()
